# -----------------------------------------------------------------------
# Applies the commit "Update racial misclassification language to neutral
# terms for non-electoral resumes" to the cartographic_professional /
# software_engineering resume document.
#
#   1. Professional summary: "affecting all Black and Asian-American
#      voters" -> "affecting 50M voters"
#   2. Siege Analytics bullet: same phrase neutralized, with "50M" bolded
#      and colored to match the other stat callouts in that bullet.
#   3. Reorders PROFESSIONAL EXPERIENCE: "Data Products Manager -
#      Helm/Murmuration" moves to right after "Partner - Siege Analytics";
#      "Research Director - PCCC" moves to right after
#      "Software Engineer - Mautinoa Technologies".
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Professional summary paragraph
# ---------------------------------------------------------------------
$sumRange = $d.Content
$null = $sumRange.Find.Execute(
    "Software engineer with 15+ years building systems that matter. Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Software engineer with 15+ years building systems that matter. Discovered systematic demographic coding errors affecting 50M voters, developed",
    2)

# ---------------------------------------------------------------------
# 2) Siege Analytics bullet - neutralize + bold/colorize "50M"
#    (scoped to that one bullet paragraph so the earlier "50M" that now
#    also appears in the professional summary is never touched)
# ---------------------------------------------------------------------
$locateBullet = $d.Content
$null = $locateBullet.Find.Execute("race coding errors affecting", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bulletParaRange = $locateBullet.Paragraphs(1).Range

$null = $bulletParaRange.Find.Execute(
    "affecting all Black and Asian-American",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M",
    2)

$boldScope = $d.Content
$null = $boldScope.Find.Execute("race coding errors affecting", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boldParaRange = $boldScope.Paragraphs(1).Range
$found50m = $boldParaRange.Find.Execute("50M", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found50m) {
    $boldParaRange.Font.Bold = $true
    $boldParaRange.Font.Color = 5258796
}

# ---------------------------------------------------------------------
# Helper: cut a 5-paragraph block (heading + 4 body paragraphs) that
# starts with $headingText, and paste it immediately before the
# paragraph that starts with $targetText. Restores the Heading 3 style
# on the moved heading paragraph (Cut/Paste only carries plain text).
# ---------------------------------------------------------------------
function Move-ExperienceBlock($headingText, $targetText) {
    $findRange = $d.Content
    $null = $findRange.Find.Execute($headingText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $startPara = $findRange.Paragraphs(1)
    $endPara = $startPara.Next(4)
    $blockStart = $startPara.Range.Start
    $blockEnd = $endPara.Range.End
    $block = $d.Range($blockStart, $blockEnd)
    $block.Cut()

    $targetRange = $d.Content
    $null = $targetRange.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $insertAt = $targetRange.Start
    $insertPoint = $d.Range($insertAt, $insertAt)
    $insertPoint.Paste()

    # Restore the Heading 3 paragraph style on the (re)moved heading line.
    $restoreRange = $d.Content
    $null = $restoreRange.Find.Execute($headingText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $restoreRange.Paragraphs(1).Style = "Heading 3"
}

# ---------------------------------------------------------------------
# 3) Move "Data Products Manager - Helm/Murmuration" block to right
#    after the "Partner - Siege Analytics" section (i.e. immediately
#    before "Software Engineer - Mautinoa Technologies").
# ---------------------------------------------------------------------
Move-ExperienceBlock "Data Products Manager - Helm/Murmuration" "Software Engineer - Mautinoa Technologies"

# Restore bold/color on the "57%" stat inside the moved block.
$locateEtl = $d.Content
$null = $locateEtl.Find.Execute("Modernized legacy ETL processes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$etlParaRange = $locateEtl.Paragraphs(1).Range
$pctFound = $etlParaRange.Find.Execute("57%", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($pctFound) {
    $etlParaRange.Font.Bold = $true
    $etlParaRange.Font.Color = 5258796
}

# ---------------------------------------------------------------------
# 4) Move "Research Director - PCCC" block to right after the
#    "Software Engineer - Mautinoa Technologies" section (i.e.
#    immediately before "Software Engineer - Salsa Labs").
# ---------------------------------------------------------------------
Move-ExperienceBlock "Research Director - PCCC" "Software Engineer - Salsa Labs"

Write-Output "done"
